# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 344
    $ws.Range("F3").Value = 90
    $ws.Range("F4").Value = 1524
    $ws.Range("F8").Value = 51
    $ws.Range("F9").Value = 342
}
